$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Remove the hyperlink that previously covered AP10:AP24 (but leave the
# AP9 hyperlink intact), since those cells are no longer independent DOI
# strings but formula-derived values.
foreach ($hl in @($ws.Hyperlinks)) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$AP$10:$AP$24') {
        $hl.Delete()
    }
}

# AP10 becomes a formula referencing AP9
$ws.Range("AP10").Formula = "=AP9"

# AP11:AP24 become a shared formula referencing the cell above
$ws.Range("AP11:AP24").Formula = "=AP10"

# Update the selection to match AP10:AP10:AP24
$ws.Range("AP10:AP24").Select()
